$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.89632680511049
$ws.Range("D2").Value = 0.9163268051104865
$ws.Range("E2").Value = 0.8396548137639914
$ws.Range("C3").Value = 29.74049686643244
$ws.Range("D3").Value = 0.5904968664324421
$ws.Range("E3").Value = 0.3486865492665334
$ws.Range("C4").Value = 30.17802841988026
$ws.Range("D4").Value = 0.8280284198802619
$ws.Range("E4").Value = 0.6856310641294032
$ws.Range("C5").Value = 29.89923241292755
$ws.Range("D5").Value = 0.529232412927545
$ws.Range("E5").Value = 0.2800869468931115
$ws.Range("C6").Value = 28.9817084540674
$ws.Range("D6").Value = -0.5582915459325974
$ws.Range("E6").Value = 0.3116894502598095
$ws.Range("C7").Value = 29.32096095580385
$ws.Range("D7").Value = -0.2290390441961527
$ws.Range("E7").Value = 0.05245888376628721
$ws.Range("C8").Value = 29.6377051150338
$ws.Range("D8").Value = -0.1122948849662038
$ws.Range("E8").Value = 0.01261014118957294
$ws.Range("C9").Value = 30.13562491976559
$ws.Range("D9").Value = 0.29562491976559
$ws.Range("E9").Value = 0.08739409318641153
$ws.Range("C10").Value = 30.02736367531402
$ws.Range("D10").Value = 0.2173636753140187
$ws.Range("E10").Value = 0.04724696734601813
$ws.Range("C11").Value = 29.79807864747706
$ws.Range("D11").Value = -0.1219213525229428
$ws.Range("E11").Value = 0.01486481620102369
$ws.Range("C12").Value = 29.73481586349194
$ws.Range("D12").Value = -0.2451841365080618
$ws.Range("E12").Value = 0.06011526079520389
$ws.Range("C13").Value = 29.74052743269053
$ws.Range("D13").Value = -0.2994725673094685
$ws.Range("E13").Value = 0.08968381857092413
$ws.Range("C14").Value = 29.81772104667583
$ws.Range("D14").Value = -0.3922789533241691
$ws.Range("E14").Value = 0.1538827772211056
$ws.Range("C15").Value = 29.83126696289785
$ws.Range("D15").Value = -0.3887330371021527
$ws.Range("E15").Value = 0.1511133741346636
$ws.Range("C16").Value = 29.97346329811299
$ws.Range("D16").Value = -0.4065367018870063
$ws.Range("E16").Value = 0.1652720899811646
$ws.Range("C17").Value = 30.04707284956853
$ws.Range("D17").Value = -0.392927150431472
$ws.Range("E17").Value = 0.1543917455461966
$ws.Range("C18").Value = 29.91660694866449
$ws.Range("D18").Value = -0.5633930513355061
$ws.Range("E18").Value = 0.3174117302931322
$ws.Range("C19").Value = 29.99936347057776
$ws.Range("D19").Value = -0.6906365294222425
$ws.Range("E19").Value = 0.4769788157724
$ws.Range("C20").Value = 30.42549840475123
$ws.Range("D20").Value = -0.3245015952487691
$ws.Range("E20").Value = 0.105301285318996
$ws.Range("C21").Value = 30.47888267616282
$ws.Range("D21").Value = -0.4611173238371791
$ws.Range("E21").Value = 0.2126291863427619
$ws.Range("C22").Value = 30.71352857353821
$ws.Range("D22").Value = -0.2364714264617938
$ws.Range("E22").Value = 0.05591873553287553
$ws.Range("C23").Value = 30.96141697364089
$ws.Range("D23").Value = -0.05858302635910917
$ws.Range("E23").Value = 0.00343197097739208
$ws.Range("C24").Value = 31.1678367404652
$ws.Range("D24").Value = 0.04783674046520048
$ws.Range("E24").Value = 0.002288353738334949
$ws.Range("C25").Value = 31.10946212854389
$ws.Range("D25").Value = -0.1705378714561085
$ws.Range("E25").Value = 0.0290831656007802
$ws.Range("C26").Value = 31.22930082036947
$ws.Range("D26").Value = -0.1506991796305321
$ws.Range("E26").Value = 0.02271024274131539
$ws.Range("C27").Value = 31.5887951343716
$ws.Range("D27").Value = 0.008795134371599289
$ws.Range("E27").Value = 0.00007735438861448722
$ws.Range("C28").Value = 31.96396303559739
$ws.Range("D28").Value = 0.3139630355973893
$ws.Range("E28").Value = 0.09857278772152751
$ws.Range("C29").Value = 32.81493599015737
$ws.Range("D29").Value = 0.9349359901573742
$ws.Range("E29").Value = 0.8741053056915496
$ws.Range("C30").Value = 32.97279510279719
$ws.Range("D30").Value = 0.6927951027971844
$ws.Range("E30").Value = 0.4799650544597612
$ws.Range("C31").Value = 33.04564124508133
$ws.Range("D31").Value = 0.5956412450813247
$ws.Range("E31").Value = 0.3547884928420307
$ws.Range("C32").Value = 33.20541006560473
$ws.Range("D32").Value = 0.3554100656047297
$ws.Range("E32").Value = 0.1263163147331582
$ws.Range("C33").Value = 33.28618459035178
$ws.Range("D33").Value = 0.3861845903517818
$ws.Range("E33").Value = 0.1491385378251735
$ws.Range("C34").Value = 33.33682847309746
$ws.Range("D34").Value = 0.2368284730974537
$ws.Range("E34").Value = 0.05608772566967135
$ws.Range("C35").Value = 33.6759927999388
$ws.Range("D35").Value = 0.2759927999387983
$ws.Range("E35").Value = 0.07617202561805755
$ws.Range("C36").Value = 33.7170788925518
$ws.Range("D36").Value = 0.01707889255179396
$ws.Range("E36").Value = 0.000291688570795723
$ws.Range("C37").Value = 33.82441943356593
$ws.Range("D37").Value = -0.2755805664340727
$ws.Range("E37").Value = 0.07594464859612435
$ws.Range("C38").Value = 34.25753728015376
$ws.Range("D38").Value = -0.142462719846236
$ws.Range("E38").Value = 0.02029562654598714
$ws.Range("C39").Value = 34.55309186796224
$ws.Range("D39").Value = -0.3469081320377612
$ws.Range("E39").Value = 0.1203452520739287
$ws.Range("C40").Value = 35.46894970214569
$ws.Range("D40").Value = 0.1689497021456958
$ws.Range("E40").Value = 0.02854400185511934
$ws.Range("C41").Value = 35.73876694052214
$ws.Range("D41").Value = 0.03876694052213736
$ws.Range("E41").Value = 0.001502875677446936
$ws.Range("C42").Value = 36.10449998253414
$ws.Range("D42").Value = -0.1955000174658537
$ws.Range("E42").Value = 0.03822025682914908
$ws.Range("C43").Value = 36.68021591248571
$ws.Range("D43").Value = -0.119784087514283
$ws.Range("E43").Value = 0.01434822762162942
$ws.Range("C44").Value = 36.80267619744119
$ws.Range("D44").Value = -0.4973238025588103
$ws.Range("E44").Value = 0.2473309645915545
$ws.Range("C45").Value = 37.74713299675034
$ws.Range("D45").Value = -0.1528670032496606
$ws.Range("E45").Value = 0.02336832068253175
$ws.Range("C46").Value = 38.41971972881053
$ws.Range("D46").Value = -0.08028027118946568
$ws.Range("E46").Value = 0.006444921942254153
$ws.Range("C47").Value = 39.09180267373289
$ws.Range("D47").Value = 0.1918026737328873
$ws.Range("E47").Value = 0.03678826565108441
$ws.Range("C48").Value = 39.59583476909648
$ws.Range("D48").Value = 0.1958347690964786
$ws.Range("E48").Value = 0.03835125678707108
$ws.Range("C49").Value = 39.71093423148177
$ws.Range("D49").Value = -0.1890657685182262
$ws.Range("E49").Value = 0.0357458648253875
$ws.Range("C50").Value = 40.30340398662099
$ws.Range("D50").Value = 0.2034039866209838
$ws.Range("E50").Value = 0.04137318177330936
$ws.Range("C51").Value = 40.58522491614487
$ws.Range("D51").Value = -0.01477508385512749
$ws.Range("E51").Value = 0.000218303102926049
$ws.Range("C52").Value = 0.2241264109621923
$ws.Range("E52").Value = 7.624873534645253
$ws.Range("E53").Value = 0.1524974706929051
